# Weekly update: two new "Ají" price records are reported for Femacal de
# La Calera (Coquimbo). They are prepended to the existing series, which
# pushes every pre-existing detail row (236-280) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right above the current first detail row of this
# block (row 236). Everything at/after row 236 shifts down to 238+.
$ws.Rows("236:237").Insert()

# --- New row 236 -----------------------------------------------------
$ws.Cells.Item(236, 1).Value = 3
$ws.Cells.Item(236, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(236, 3).Value = "Coquimbo"
$ws.Cells.Item(236, 4).Value = 44476
$ws.Cells.Item(236, 5).Value = 5
$ws.Cells.Item(236, 6).Value = 100112021
$ws.Cells.Item(236, 7).Value = "Ají"
$ws.Cells.Item(236, 8).Value = "Americana (o)"
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 65
$ws.Cells.Item(236, 11).Value = 40000
$ws.Cells.Item(236, 12).Value = 41000
$ws.Cells.Item(236, 13).Value = 40462
$ws.Cells.Item(236, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(236, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(236, 16).Value = 2697
$ws.Cells.Item(236, 17).Value = 15
$ws.Cells.Item(236, 18).Value = "Hortaliza"

# --- New row 237 -----------------------------------------------------
$ws.Cells.Item(237, 1).Value = 3
$ws.Cells.Item(237, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(237, 3).Value = "Coquimbo"
$ws.Cells.Item(237, 4).Value = 44476
$ws.Cells.Item(237, 5).Value = 5
$ws.Cells.Item(237, 6).Value = 100112021
$ws.Cells.Item(237, 7).Value = "Ají"
$ws.Cells.Item(237, 8).Value = "Americana (o)"
$ws.Cells.Item(237, 9).Value = "Primera"
$ws.Cells.Item(237, 10).Value = 35
$ws.Cells.Item(237, 11).Value = 85000
$ws.Cells.Item(237, 12).Value = 85000
$ws.Cells.Item(237, 13).Value = 85000
$ws.Cells.Item(237, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(237, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(237, 16).Value = 3400
$ws.Cells.Item(237, 17).Value = 25
$ws.Cells.Item(237, 18).Value = "Hortaliza"
